$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the data row for the new snapshot date (7 March 2022 update -> 03-Feb-22 entry)
$ws.Range("A8").Value = 44595
$ws.Range("B8").Value = 0.04
$ws.Range("C8").Value = 0.19
$ws.Range("D8").Value = 0.12
$ws.Range("E8").Value = 0.5
$ws.Range("F8").Value = 0.08

# Update view state: scroll/selection
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("F8").Select()
